$wb = $excel.ActiveWorkbook

$wsSite = $wb.Worksheets.Item("site_metrics")
$wsDur = $wb.Worksheets.Item("mk_duration")
$wsIntra = $wb.Worksheets.Item("mk_intra_annual")

# --- site_metrics ---
$wsSite.Range("O2").Value = 0.02549195762010053
$wsSite.Range("N5").Value = 8.207563025210083
$wsSite.Range("O7").Value = 0.004813189684878138
$wsSite.Range("AK13").Value = $true
$wsSite.Range("AK14").Value = $true
$wsSite.Range("O18").Value = 0.009845425644055645
$wsSite.Range("O30").Value = 0.07811873031498129
$wsSite.Range("N32").Value = 3.103713151927437
$wsSite.Range("O33").Value = 0.1402879772438876
$wsSite.Range("O34").Value = 0.01860710227192392
$wsSite.Range("AK36").Value = $true
$wsSite.Range("AK37").Value = $true
$wsSite.Range("O39").Value = 0.412524850340755
$wsSite.Range("AK40").Value = $true
$wsSite.Range("AK41").Value = $true
$wsSite.Range("AK43").Value = $true
$wsSite.Range("AK44").Value = $true
$wsSite.Range("AK52").Value = $true
$wsSite.Range("O56").Value = 0.08052547057925852
$wsSite.Range("O57").Value = 0.08906252763741357
$wsSite.Range("AK60").Value = $true
$wsSite.Range("N61").Value = 3.469658995974786
$wsSite.Range("N65").Value = 3.529120879120879
$wsSite.Range("O67").Value = 0.01099562600072403
$wsSite.Range("AK69").Value = $true
$wsSite.Range("N70").Value = 10.93350694444444
$wsSite.Range("O70").Value = 0.008225483830074253
$wsSite.Range("Q70").Value = 2.6875
$wsSite.Range("O71").Value = 0.003885208861282886
$wsSite.Range("AK73").Value = $true
$wsSite.Range("N75").Value = 3.402378983024145
$wsSite.Range("O75").Value = 0.001079330339997736
$wsSite.Range("N78").Value = 3.49597162097162
$wsSite.Range("O78").Value = 0.006176649759310394
$wsSite.Range("O80").Value = 0.004193311531737252

# --- mk_duration ---
$wsDur.Range("K5").Value = "no trend"
$wsDur.Range("L5").Value = $false
$wsDur.Range("M5").Value = 0.09344799392558367
$wsDur.Range("N5").Value = -1.677483394552118
$wsDur.Range("O5").Value = -0.2142857142857143
$wsDur.Range("P5").Value = -87
$wsDur.Range("Q5").Value = 2628.333333333333
$wsDur.Range("R5").Value = -0.1091269841269841
$wsDur.Range("S5").Value = 3.527777777777778
$wsDur.Range("K26").Value = "no trend"
$wsDur.Range("L26").Value = $false
$wsDur.Range("M26").Value = 0.1099636292968564
$wsDur.Range("N26").Value = 1.598356637186201
$wsDur.Range("O26").Value = 0.196969696969697
$wsDur.Range("P26").Value = 104
$wsDur.Range("Q26").Value = 4152.666666666667
$wsDur.Range("R26").Value = 0.007905982905982906
$wsDur.Range("S26").Value = 1.54017094017094
$wsDur.Range("M31").Value = 0.9523825905788506
$wsDur.Range("N31").Value = -0.0597150429446484
$wsDur.Range("O31").Value = -0.006205673758865249
$wsDur.Range("P31").Value = -7
$wsDur.Range("Q31").Value = 10095.66666666667
$wsDur.Range("K32").Value = "no trend"
$wsDur.Range("L32").Value = $false
$wsDur.Range("M32").Value = 0.9401083689080136
$wsDur.Range("N32").Value = 0.07513365721922702
$wsDur.Range("O32").Value = 0.01231527093596059
$wsDur.Range("P32").Value = 5
$wsDur.Range("Q32").Value = 2834.333333333333
$wsDur.Range("R32").Value = 0
$wsDur.Range("S32").Value = 2.714285714285714
$wsDur.Range("K39").Value = "no trend"
$wsDur.Range("L39").Value = $false
$wsDur.Range("M39").Value = 0.6157503017697474
$wsDur.Range("N39").Value = 0.5018823210851089
$wsDur.Range("O39").Value = 0.05666666666666666
$wsDur.Range("P39").Value = 17
$wsDur.Range("Q39").Value = 1016.333333333333
$wsDur.Range("M61").Value = 0.7214746342177878
$wsDur.Range("N61").Value = 0.3564886717836127
$wsDur.Range("O61").Value = 0.03875968992248062
$wsDur.Range("P61").Value = 35
$wsDur.Range("Q61").Value = 9096.333333333334
$wsDur.Range("R61").Value = 0.007575757575757569
$wsDur.Range("S61").Value = 2.440909090909091
$wsDur.Range("M65").Value = 0.8123923688472134
$wsDur.Range("N65").Value = -0.2373408560833091
$wsDur.Range("O65").Value = -0.03439153439153439
$wsDur.Range("P65").Value = -13
$wsDur.Range("Q65").Value = 2556.333333333333
$wsDur.Range("R65").Value = -0.03003968253968263
$wsDur.Range("S65").Value = 3.755535714285716
$wsDur.Range("M66").Value = 0.3660827184679289
$wsDur.Range("N66").Value = -0.9038353419323745
$wsDur.Range("O66").Value = -0.1051051051051051
$wsDur.Range("P66").Value = -70
$wsDur.Range("Q66").Value = 5828
$wsDur.Range("R66").Value = -0.03858560794044666
$wsDur.Range("S66").Value = 3.617617866004963
$wsDur.Range("D70").Value = 0.6339498028447381
$wsDur.Range("E70").Value = -0.476174867647586
$wsDur.Range("F70").Value = -0.06236559139784946
$wsDur.Range("G70").Value = -29
$wsDur.Range("I70").Value = -0.05128205128205127
$wsDur.Range("J70").Value = 8.435897435897436
$wsDur.Range("M70").Value = 0.2738534999027467
$wsDur.Range("N70").Value = -1.094231407679569
$wsDur.Range("O70").Value = -0.1044897959183674
$wsDur.Range("P70").Value = -128
$wsDur.Range("M71").Value = 0.4888249449426469
$wsDur.Range("N71").Value = -0.6921789680474556
$wsDur.Range("O71").Value = -0.07149758454106281
$wsDur.Range("P71").Value = -74
$wsDur.Range("Q71").Value = 11122.66666666667
$wsDur.Range("R71").Value = -0.01388888888888889
$wsDur.Range("S71").Value = 2.3125
$wsDur.Range("M75").Value = 0.5748247878461354
$wsDur.Range("N75").Value = -0.5609600264316629
$wsDur.Range("O75").Value = -0.06666666666666667
$wsDur.Range("P75").Value = -42
$wsDur.Range("Q75").Value = 5342
$wsDur.Range("R75").Value = -0.01044383270549599
$wsDur.Range("S75").Value = 2.18276707234618
$wsDur.Range("K76").Value = "no trend"
$wsDur.Range("L76").Value = $false
$wsDur.Range("M76").Value = 0.5358598692967174
$wsDur.Range("N76").Value = 0.619085751600795
$wsDur.Range("O76").Value = 0.06829268292682927
$wsDur.Range("P76").Value = 56
$wsDur.Range("Q76").Value = 7892.666666666667
$wsDur.Range("R76").Value = 0.008130215649012631
$wsDur.Range("S76").Value = 1.837395687019747
$wsDur.Range("M78").Value = 0.8869230368864609
$wsDur.Range("N78").Value = 0.142198727497721
$wsDur.Range("O78").Value = 0.01545893719806763
$wsDur.Range("P78").Value = 16
$wsDur.Range("Q78").Value = 11127.33333333333
$wsDur.Range("S78").Value = 3.171428571428572
$wsDur.Range("K80").Value = "no trend"
$wsDur.Range("L80").Value = $false
$wsDur.Range("M80").Value = 0.3069912409644038
$wsDur.Range("N80").Value = 1.021555684791887
$wsDur.Range("O80").Value = 0.1333333333333333
$wsDur.Range("P80").Value = 58
$wsDur.Range("Q80").Value = 3113.333333333333
$wsDur.Range("R80").Value = 0.02500000000000002
$wsDur.Range("S80").Value = 2.6375
$wsDur.Range("M82").Value = 0.858673521284679
$wsDur.Range("N82").Value = 0.1780629720586193
$wsDur.Range("O82").Value = 0.0231729055258467
$wsDur.Range("P82").Value = 13
$wsDur.Range("Q82").Value = 4541.666666666667
$wsDur.Range("R82").Value = 0.003654970760233928
$wsDur.Range("S82").Value = 2.38969298245614

# --- mk_intra_annual ---
$wsIntra.Range("K5").Value = "no trend"
$wsIntra.Range("S5").Value = 1
$wsIntra.Range("M5").Value = 0.410746198342077
$wsIntra.Range("N5").Value = -0.8225811910567383
$wsIntra.Range("O5").Value = -0.1059113300492611
$wsIntra.Range("P5").Value = -43
$wsIntra.Range("Q5").Value = 2607
$wsIntra.Range("L5").Value = $false
$wsIntra.Range("M26").Value = 0.3249311654768372
$wsIntra.Range("N26").Value = -0.984375
$wsIntra.Range("O26").Value = -0.1212121212121212
$wsIntra.Range("P26").Value = -64
$wsIntra.Range("Q26").Value = 4096
$wsIntra.Range("R26").Value = -0.08514492753623187
$wsIntra.Range("S26").Value = 12.36231884057971
$wsIntra.Range("M31").Value = 0.8101656932688071
$wsIntra.Range("N31").Value = -0.2402122813141875
$wsIntra.Range("O31").Value = -0.02216312056737589
$wsIntra.Range("P31").Value = -25
$wsIntra.Range("Q31").Value = 9982.333333333334
$wsIntra.Range("K32").Value = "no trend"
$wsIntra.Range("L32").Value = $false
$wsIntra.Range("M32").Value = 0.9849455743562765
$wsIntra.Range("N32").Value = -0.01886904412032354
$wsIntra.Range("O32").Value = -0.004926108374384237
$wsIntra.Range("P32").Value = -2
$wsIntra.Range("Q32").Value = 2808.666666666667
$wsIntra.Range("R32").Value = 0
$wsIntra.Range("S32").Value = 6
$wsIntra.Range("K39").Value = "no trend"
$wsIntra.Range("L39").Value = $false
$wsIntra.Range("M39").Value = 0.6128224102281419
$wsIntra.Range("N39").Value = 0.5060487741189207
$wsIntra.Range("O39").Value = 0.05666666666666666
$wsIntra.Range("P39").Value = 17
$wsIntra.Range("Q39").Value = 999.6666666666666
$wsIntra.Range("M61").Value = 0.3002074359872218
$wsIntra.Range("N61").Value = -1.035988653519739
$wsIntra.Range("O61").Value = -0.1096345514950166
$wsIntra.Range("P61").Value = -99
$wsIntra.Range("Q61").Value = 8948.333333333334
$wsIntra.Range("R61").Value = -0.03571428571428571
$wsIntra.Range("S61").Value = 5.75
$wsIntra.Range("M65").Value = 0.6755034011772314
$wsIntra.Range("N65").Value = -0.4186069613366101
$wsIntra.Range("O65").Value = -0.0582010582010582
$wsIntra.Range("P65").Value = -22
$wsIntra.Range("Q65").Value = 2516.666666666667
$wsIntra.Range("S65").Value = 4.5
$wsIntra.Range("M66").Value = 0.6251697703125925
$wsIntra.Range("N66").Value = -0.4885366530433578
$wsIntra.Range("O66").Value = -0.05705705705705705
$wsIntra.Range("P66").Value = -38
$wsIntra.Range("Q66").Value = 5736
$wsIntra.Range("S66").Value = 4
$wsIntra.Range("M71").Value = 0.03614983237685898
$wsIntra.Range("N71").Value = 2.095238095238095
$wsIntra.Range("O71").Value = 0.2135265700483092
$wsIntra.Range("P71").Value = 221
$wsIntra.Range("Q71").Value = 11025
$wsIntra.Range("R71").Value = 0.1
$wsIntra.Range("S71").Value = 3.75
$wsIntra.Range("M75").Value = 0.923614755097441
$wsIntra.Range("N75").Value = -0.09588141518882391
$wsIntra.Range("O75").Value = -0.0126984126984127
$wsIntra.Range("P75").Value = -8
$wsIntra.Range("Q75").Value = 5330
$wsIntra.Range("S75").Value = 5
$wsIntra.Range("K76").Value = "no trend"
$wsIntra.Range("L76").Value = $false
$wsIntra.Range("M76").Value = 0.4041663847484855
$wsIntra.Range("N76").Value = 0.8342033836521727
$wsIntra.Range("O76").Value = 0.09146341463414634
$wsIntra.Range("P76").Value = 75
$wsIntra.Range("Q76").Value = 7869
$wsIntra.Range("R76").Value = 0.06559139784946236
$wsIntra.Range("S76").Value = 5.688172043010753
$wsIntra.Range("M78").Value = 0.4913828936852183
$wsIntra.Range("N78").Value = -0.6881109845448636
$wsIntra.Range("O78").Value = -0.07053140096618357
$wsIntra.Range("P78").Value = -73
$wsIntra.Range("Q78").Value = 10948.33333333333
$wsIntra.Range("M80").Value = 0.3146467998598692
$wsIntra.Range("N80").Value = -1.005519426749776
$wsIntra.Range("O80").Value = -0.1310344827586207
$wsIntra.Range("P80").Value = -57
$wsIntra.Range("Q80").Value = 3101.666666666667
$wsIntra.Range("R80").Value = -0.08333333333333333
$wsIntra.Range("S80").Value = 6.708333333333333
$wsIntra.Range("M82").Value = 0.730547120402357
$wsIntra.Range("N82").Value = -0.3443978302354009
$wsIntra.Range("O82").Value = -0.0427807486631016
$wsIntra.Range("Q82").Value = 4460
$wsIntra.Range("S82").Value = 5
